# Update EC database and add part 1 of new account statement period
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update "VALOR MORA" total amount
$ws.Range("E11").Value = 227760

# 2) Update worker/period counts
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# 3) Correct worker identity on row 17 (was 73214199 / OSWALDO SIERRA SANTIAGO)
$ws.Range("C17").Value = "1047461168"
$ws.Range("D17").Value = "ALEXANDER ORTEGA VASQUEZ"

# 4) Insert a new data row before row 18, copying row 17's formatting, for the new period
$ws.Rows.Item(17).Copy()
$ws.Rows.Item(18).Insert()

# Ensure full-grid thin borders on the newly inserted row match the table style
$rng = $ws.Range("B18:J18")
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(8).LineStyle = 1
$rng.Borders.Item(9).LineStyle = 1
$rng.Borders.Item(10).LineStyle = 1
$rng.Borders.Item(11).LineStyle = 1
$rng.Borders.Item(12).LineStyle = 1
$ws.Range("B18").Borders.Item(7).ColorIndex = 1
$ws.Range("J18").Borders.Item(10).ColorIndex = 1

# New row 18: Rafael Alfonso Ortega Vasquez, period 2508
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73199389"
$ws.Range("D18").Value = "RAFAEL ALFONSO ORTEGA VASQUEZ"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 (previously row 18, shifted down by the insert): Alexander Ortega Vasquez, period 2508
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047461168"
$ws.Range("D19").Value = "ALEXANDER ORTEGA VASQUEZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

Write-Host "Done"
